$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 111112390
$ws.Range("I127").Value = 1333.3334
$ws.Range("J127").Value = 166667920
$ws.Range("K127").Value = 4000.0002
$ws.Range("L127").Value = 500003760
$ws.Range("M127").Value = 959.9998000000001
$ws.Range("N127").Value = -500013680
$ws.Range("H138").Value = 3599.6365
$ws.Range("J138").Value = 4475.305
$ws.Range("L138").Value = 13425.915
$ws.Range("N138").Value = -23705.915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1578.7106
$ws.Range("I2").Value = 1880.9546
$ws.Range("J2").Value = 1163.125
$ws.Range("K2").Value = 1880.9546
$ws.Range("L2").Value = 1163.125
$ws.Range("M2").Value = -1767.9546
$ws.Range("N2").Value = -1389.125
$ws.Range("H32").Value = 8016.375
$ws.Range("I32").Value = 7326.1616
$ws.Range("K32").Value = 7326.1616
$ws.Range("M32").Value = -7039.1616
$ws.Range("H74").Value = 109576.7
$ws.Range("I74").Value = 120282.2
$ws.Range("J74").Value = 29285.5
$ws.Range("K74").Value = 120282.2
$ws.Range("L74").Value = 29285.5
$ws.Range("M74").Value = -119408.2
$ws.Range("N74").Value = -31033.5
$ws.Range("H77").Value = 109576.7
$ws.Range("I77").Value = 120282.2
$ws.Range("J77").Value = 29285.5
$ws.Range("K77").Value = 601411
$ws.Range("L77").Value = 146427.5
$ws.Range("M77").Value = -597043
$ws.Range("N77").Value = -155163.5
$ws.Range("H110").Value = 1951
$ws.Range("I110").Value = 1793.3846
$ws.Range("J110").Value = 4000
$ws.Range("K110").Value = 1793.3846
$ws.Range("L110").Value = 4000
$ws.Range("M110").Value = 251.6153999999999
$ws.Range("N110").Value = -8090
$ws.Range("H116").Value = 1578.7106
$ws.Range("I116").Value = 1880.9546
$ws.Range("J116").Value = 1163.125
$ws.Range("K116").Value = 1880.9546
$ws.Range("L116").Value = 1163.125
$ws.Range("M116").Value = 413.0454
$ws.Range("N116").Value = -5751.125
$ws.Range("H125").Value = 65485.91
$ws.Range("J125").Value = 65485.91
$ws.Range("L125").Value = 65485.91
$ws.Range("N125").Value = -75325.91

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1578.7106
$ws.Range("I3").Value = 1880.9546
$ws.Range("J3").Value = 1163.125
$ws.Range("K3").Value = 1880.9546
$ws.Range("L3").Value = 1163.125
$ws.Range("M3").Value = -1766.9546
$ws.Range("N3").Value = -1391.125
$ws.Range("H126").Value = 80780
$ws.Range("J126").Value = 80780
$ws.Range("L126").Value = 80780
$ws.Range("N126").Value = -90660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3899.0378
$ws.Range("I31").Value = 2091.6099
$ws.Range("K31").Value = 2091.6099
$ws.Range("M31").Value = -1796.6099
$ws.Range("H34").Value = 3899.0378
$ws.Range("I34").Value = 2091.6099
$ws.Range("K34").Value = 2091.6099
$ws.Range("M34").Value = -1889.6099
$ws.Range("H41").Value = 9999.5
$ws.Range("I41").Value = 9999.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 9999.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -9571.5
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 24394.666
$ws.Range("J50").Value = 24394.666
$ws.Range("L50").Value = 24394.666
$ws.Range("N50").Value = -25644.666
$ws.Range("H51").Value = 28000
$ws.Range("J51").Value = 28000
$ws.Range("L51").Value = 28000
$ws.Range("N51").Value = -29472
$ws.Range("H60").Value = 13666.667
$ws.Range("I60").Value = 13666.667
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 13666.667
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -13155.667
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 28000
$ws.Range("J61").Value = 28000
$ws.Range("L61").Value = 28000
$ws.Range("N61").Value = -28696
$ws.Range("H105").Value = 799.6667
$ws.Range("J105").Value = 1500
$ws.Range("L105").Value = 1500
$ws.Range("N105").Value = -4994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2448.5
$ws.Range("J75").Value = 2448.5
$ws.Range("L75").Value = 7345.5
$ws.Range("N75").Value = -9341.5
$ws.Range("H78").Value = 2448.5
$ws.Range("J78").Value = 2448.5
$ws.Range("L78").Value = 22036.5
$ws.Range("N78").Value = -32020.5
$ws.Range("H86").Value = 1999
$ws.Range("I86").Value = 2198.8
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 6596.400000000001
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -5410.400000000001
$ws.Range("N86").Value = -5372
$ws.Range("H87").Value = 7111.933
$ws.Range("I87").Value = 2171.3333
$ws.Range("J87").Value = 8347.083000000001
$ws.Range("K87").Value = 6513.999899999999
$ws.Range("L87").Value = 25041.249
$ws.Range("M87").Value = -5265.999899999999
$ws.Range("N87").Value = -27537.249
$ws.Range("H89").Value = 1999
$ws.Range("I89").Value = 2198.8
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 19789.2
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -13861.2
$ws.Range("N89").Value = -20856
$ws.Range("H90").Value = 7111.933
$ws.Range("I90").Value = 2171.3333
$ws.Range("J90").Value = 8347.083000000001
$ws.Range("K90").Value = 19541.9997
$ws.Range("L90").Value = 75123.747
$ws.Range("M90").Value = -13301.9997
$ws.Range("N90").Value = -87603.747
$ws.Range("H93").Value = 4632.2856
$ws.Range("J93").Value = 4987.6665
$ws.Range("L93").Value = 14962.9995
$ws.Range("N93").Value = -18706.9995
$ws.Range("H114").Value = 460.83334
$ws.Range("I114").Value = 460.83334
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 1382.50002
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 1871.49998
$ws.Range("N114").ClearContents()
$ws.Range("H129").Value = 1790
$ws.Range("I129").Value = 2185.6667
$ws.Range("J129").Value = 1394.3334
$ws.Range("K129").Value = 6557.000100000001
$ws.Range("L129").Value = 4183.0002
$ws.Range("M129").Value = -1557.000100000001
$ws.Range("N129").Value = -14183.0002
$ws.Range("H139").Value = 1640541
$ws.Range("I139").Value = 2819491.5
$ws.Range("J139").Value = 3109.8333
$ws.Range("K139").Value = 8458474.5
$ws.Range("L139").Value = 9329.499899999999
$ws.Range("M139").Value = -8453334.5
$ws.Range("N139").Value = -19609.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1076
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1076
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1076
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1452
$ws.Range("H82").Value = 2388.5
$ws.Range("I82").Value = 1265.8334
$ws.Range("J82").Value = 4633.8335
$ws.Range("K82").Value = 1265.8334
$ws.Range("L82").Value = 4633.8335
$ws.Range("M82").Value = -904.8334
$ws.Range("N82").Value = -5355.8335
$ws.Range("H85").Value = 2388.5
$ws.Range("I85").Value = 1265.8334
$ws.Range("J85").Value = 4633.8335
$ws.Range("K85").Value = 1265.8334
$ws.Range("L85").Value = 4633.8335
$ws.Range("M85").Value = -17.83339999999998
$ws.Range("N85").Value = -7129.8335
$ws.Range("H141").Value = 53492.5
$ws.Range("J141").Value = 53492.5
$ws.Range("L141").Value = 53492.5
$ws.Range("N141").Value = -63852.5
